# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.500.46'
$ws.Range('E2').Value = '  +4.30%  '
$ws.Range('D3').Value = '4.042.95'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'519.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').Value = "'147.16"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('D7').Value = "'0.725"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +18.36%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = "'0.759"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.61%  '
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').Value = "'0.0000328"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('D12').Value = "'47.10"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.69%  '
$ws.Range('D13').Value = "'10.91"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.33%  '
$ws.Range('D14').Value = '4.695.32'
$ws.Range('E14').Value = '  +3.55%  '
$ws.Range('D15').Value = '4.048.52'
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').Value = "'21.10"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.44%  '
$ws.Range('D17').Value = "'14.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').Value = '72.337.48'
$ws.Range('E20').Value = '  +4.19%  '
$ws.Range('D21').Value = "'444.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.66%  '
$ws.Range('D22').Value = "'104.79"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.36%  '
$ws.Range('E23').Value = '  +5.96%  '
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = "'11.02"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('D28').Value = "'37.60"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.18%  '
$ws.Range('D29').Value = "'5.81"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('E30').Value = '  +10.55%  '
$ws.Range('D31').Value = "'13.61"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.45%  '
$ws.Range('E32').Value = '  +2.80%  '
$ws.Range('D33').Value = "'677.19"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('D34').Value = "'6.84"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.90%  '
$ws.Range('D35').Value = "'67.53"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('D36').Value = "'42.91"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.97%  '
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('B38').Value = 'ThetaToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D38').Value = "'3.62"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.68%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0860'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').Value = "'0.998"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E42').Value = '  +3.10%  '
$ws.Range('D43').Value = "'0.999"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = "'0.160"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.51%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = "'3.23"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.77%  '
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').Value = "'3.43"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('D49').Value = "'9.03"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.46%  '
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = "'0.000267"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.47%  '
